$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 5226042   # B2
$ws.Cells.Item(2, 5).Value = "Zaglebie Lubin II"   # E2
$ws.Cells.Item(2, 6).Value = "KKS 1925 Kalisz"   # F2
$ws.Cells.Item(2, 7).Value = 4   # G2
$ws.Cells.Item(2, 8).Value = 3   # H2
$ws.Cells.Item(2, 10).Value = 1   # J2
$ws.Cells.Item(2, 11).Value = "H"   # K2
$ws.Cells.Item(2, 12).Value = 3.1   # L2
$ws.Cells.Item(2, 13).Value = 3.6   # M2
$ws.Cells.Item(2, 14).Value = 1.95   # N2
$ws.Cells.Item(2, 15).Value = 3.3   # O2
$ws.Cells.Item(2, 16).Value = 3.6   # P2
$ws.Cells.Item(2, 17).Value = 1.85   # Q2
$ws.Cells.Item(2, 18).Value = 0.5   # R2
$ws.Cells.Item(2, 19).Value = 1.9   # S2
$ws.Cells.Item(2, 20).Value = 1.9   # T2
$ws.Cells.Item(2, 21).Value = 2.75   # U2
$ws.Cells.Item(2, 22).Value = 1.975   # V2
$ws.Cells.Item(2, 23).Value = 1.825   # W2
$ws.Cells.Item(2, 24).Value = 2.3   # X2
$ws.Cells.Item(2, 25).Value = -1   # Y2
$ws.Cells.Item(2, 27).Value = 0.8999999999999999   # AA2
$ws.Cells.Item(2, 28).Value = -1   # AB2
$ws.Cells.Item(2, 29).Value = 0.9750000000000001   # AC2
$ws.Cells.Item(2, 30).Value = -1   # AD2

# Row 3
$ws.Cells.Item(3, 2).Value = 5228058   # B3
$ws.Cells.Item(3, 5).Value = "Polonia Warsaw"   # E3
$ws.Cells.Item(3, 6).Value = "Motor Lublin"   # F3
$ws.Cells.Item(3, 7).Value = 1   # G3
$ws.Cells.Item(3, 8).Value = 1   # H3
$ws.Cells.Item(3, 10).Value = 0   # J3
$ws.Cells.Item(3, 11).Value = "D"   # K3
$ws.Cells.Item(3, 12).Value = 3   # L3
$ws.Cells.Item(3, 13).Value = 3.3   # M3
$ws.Cells.Item(3, 14).Value = 2.1   # N3
$ws.Cells.Item(3, 15).Value = 2.8   # O3
$ws.Cells.Item(3, 16).Value = 3.25   # P3
$ws.Cells.Item(3, 17).Value = 2.2   # Q3
$ws.Cells.Item(3, 18).Value = 0.25   # R3
$ws.Cells.Item(3, 19).Value = 1.8   # S3
$ws.Cells.Item(3, 20).Value = 2   # T3
$ws.Cells.Item(3, 21).Value = 2.5   # U3
$ws.Cells.Item(3, 22).Value = 1.925   # V3
$ws.Cells.Item(3, 23).Value = 1.875   # W3
$ws.Cells.Item(3, 24).Value = -1   # X3
$ws.Cells.Item(3, 25).Value = 2.25   # Y3
$ws.Cells.Item(3, 27).Value = 0.4   # AA3
$ws.Cells.Item(3, 28).Value = -0.5   # AB3
$ws.Cells.Item(3, 29).Value = -1   # AC3
$ws.Cells.Item(3, 30).Value = 0.875   # AD3

# Row 18
$ws.Cells.Item(18, 2).Value = 6808905   # B18
$ws.Cells.Item(18, 5).Value = "Olimpia Grudziadz"   # E18
$ws.Cells.Item(18, 6).Value = "Sandecja Nowy Sacz"   # F18
$ws.Cells.Item(18, 7).Value = 3   # G18
$ws.Cells.Item(18, 8).Value = 1   # H18
$ws.Cells.Item(18, 9).Value = 2   # I18
$ws.Cells.Item(18, 10).Value = 1   # J18
$ws.Cells.Item(18, 11).Value = "H"   # K18
$ws.Cells.Item(18, 12).Value = 2.5   # L18
$ws.Cells.Item(18, 13).Value = 3.2   # M18
$ws.Cells.Item(18, 14).Value = 2.5   # N18
$ws.Cells.Item(18, 15).Value = 2.15   # O18
$ws.Cells.Item(18, 16).Value = 3.25   # P18
$ws.Cells.Item(18, 17).Value = 3.1   # Q18
$ws.Cells.Item(18, 18).Value = -0.25   # R18
$ws.Cells.Item(18, 19).Value = 1.9   # S18
$ws.Cells.Item(18, 20).Value = 1.9   # T18
$ws.Cells.Item(18, 21).Value = 2.75   # U18
$ws.Cells.Item(18, 22).Value = 2   # V18
$ws.Cells.Item(18, 23).Value = 1.8   # W18
$ws.Cells.Item(18, 24).Value = 1.15   # X18
$ws.Cells.Item(18, 25).Value = -1   # Y18
$ws.Cells.Item(18, 27).Value = 0.8999999999999999   # AA18
$ws.Cells.Item(18, 28).Value = -1   # AB18
$ws.Cells.Item(18, 29).Value = 1   # AC18
$ws.Cells.Item(18, 30).Value = -1   # AD18

# Row 20
$ws.Cells.Item(20, 2).Value = 6808604   # B20
$ws.Cells.Item(20, 5).Value = "GKS Jastrzebie"   # E20
$ws.Cells.Item(20, 6).Value = "Stal Stalowa Wola"   # F20
$ws.Cells.Item(20, 7).Value = 0   # G20
$ws.Cells.Item(20, 8).Value = 0   # H20
$ws.Cells.Item(20, 9).Value = 0   # I20
$ws.Cells.Item(20, 10).Value = 0   # J20
$ws.Cells.Item(20, 11).Value = "D"   # K20
$ws.Cells.Item(20, 12).Value = 2.2   # L20
$ws.Cells.Item(20, 13).Value = 3.1   # M20
$ws.Cells.Item(20, 14).Value = 3   # N20
$ws.Cells.Item(20, 15).Value = 2.6   # O20
$ws.Cells.Item(20, 16).Value = 3.1   # P20
$ws.Cells.Item(20, 17).Value = 2.5   # Q20
$ws.Cells.Item(20, 18).Value = 0   # R20
$ws.Cells.Item(20, 19).Value = 1.95   # S20
$ws.Cells.Item(20, 20).Value = 1.85   # T20
$ws.Cells.Item(20, 21).Value = 2.5   # U20
$ws.Cells.Item(20, 22).Value = 1.85   # V20
$ws.Cells.Item(20, 23).Value = 1.95   # W20
$ws.Cells.Item(20, 24).Value = -1   # X20
$ws.Cells.Item(20, 25).Value = 2.1   # Y20
$ws.Cells.Item(20, 27).Value = 0   # AA20
$ws.Cells.Item(20, 28).Value = 0   # AB20
$ws.Cells.Item(20, 29).Value = -1   # AC20
$ws.Cells.Item(20, 30).Value = 0.95   # AD20

# Row 44
$ws.Cells.Item(44, 2).Value = 6808166   # B44
$ws.Cells.Item(44, 5).Value = "Radunia Stezyca"   # E44
$ws.Cells.Item(44, 6).Value = "Polonia Bytom"   # F44
$ws.Cells.Item(44, 8).Value = 1   # H44
$ws.Cells.Item(44, 11).Value = "D"   # K44
$ws.Cells.Item(44, 12).Value = 2.2   # L44
$ws.Cells.Item(44, 13).Value = 3.25   # M44
$ws.Cells.Item(44, 14).Value = 2.8   # N44
$ws.Cells.Item(44, 17).Value = 2.8   # Q44
$ws.Cells.Item(44, 19).Value = 2   # S44
$ws.Cells.Item(44, 20).Value = 1.8   # T44
$ws.Cells.Item(44, 21).Value = 2.25   # U44
$ws.Cells.Item(44, 22).Value = 1.875   # V44
$ws.Cells.Item(44, 23).Value = 1.925   # W44
$ws.Cells.Item(44, 24).Value = -1   # X44
$ws.Cells.Item(44, 25).Value = 2.25   # Y44
$ws.Cells.Item(44, 27).Value = -0.5   # AA44
$ws.Cells.Item(44, 28).Value = 0.4   # AB44
$ws.Cells.Item(44, 29).Value = -0.5   # AC44
$ws.Cells.Item(44, 30).Value = 0.4625   # AD44

# Row 45
$ws.Cells.Item(45, 2).Value = 6808617   # B45
$ws.Cells.Item(45, 5).Value = "GKS Jastrzebie"   # E45
$ws.Cells.Item(45, 6).Value = "Hutnik Krakow"   # F45
$ws.Cells.Item(45, 8).Value = 0   # H45
$ws.Cells.Item(45, 9).Value = 0   # I45
$ws.Cells.Item(45, 10).Value = 0   # J45
$ws.Cells.Item(45, 11).Value = "H"   # K45
$ws.Cells.Item(45, 12).Value = 2.4   # L45
$ws.Cells.Item(45, 13).Value = 3.2   # M45
$ws.Cells.Item(45, 14).Value = 2.6   # N45
$ws.Cells.Item(45, 15).Value = 2.2   # O45
$ws.Cells.Item(45, 16).Value = 3.25   # P45
$ws.Cells.Item(45, 17).Value = 2.875   # Q45
$ws.Cells.Item(45, 18).Value = -0.25   # R45
$ws.Cells.Item(45, 19).Value = 1.975   # S45
$ws.Cells.Item(45, 20).Value = 1.825   # T45
$ws.Cells.Item(45, 21).Value = 2.5   # U45
$ws.Cells.Item(45, 22).Value = 1.925   # V45
$ws.Cells.Item(45, 23).Value = 1.875   # W45
$ws.Cells.Item(45, 24).Value = 1.2   # X45
$ws.Cells.Item(45, 26).Value = -1   # Z45
$ws.Cells.Item(45, 27).Value = 0.9750000000000001   # AA45
$ws.Cells.Item(45, 28).Value = -1   # AB45
$ws.Cells.Item(45, 29).Value = -1   # AC45
$ws.Cells.Item(45, 30).Value = 0.875   # AD45

# Row 46
$ws.Cells.Item(46, 2).Value = 6808165   # B46
$ws.Cells.Item(46, 5).Value = "Olimpia Grudziadz"   # E46
$ws.Cells.Item(46, 6).Value = "Lech Poznan II"   # F46
$ws.Cells.Item(46, 8).Value = 2   # H46
$ws.Cells.Item(46, 9).Value = 1   # I46
$ws.Cells.Item(46, 10).Value = 1   # J46
$ws.Cells.Item(46, 11).Value = "A"   # K46
$ws.Cells.Item(46, 12).Value = 1.909   # L46
$ws.Cells.Item(46, 14).Value = 3.5   # N46
$ws.Cells.Item(46, 15).Value = 1.615   # O46
$ws.Cells.Item(46, 16).Value = 3.6   # P46
$ws.Cells.Item(46, 17).Value = 4.5   # Q46
$ws.Cells.Item(46, 18).Value = -0.75   # R46
$ws.Cells.Item(46, 19).Value = 1.875   # S46
$ws.Cells.Item(46, 20).Value = 1.925   # T46
$ws.Cells.Item(46, 21).Value = 2.75   # U46
$ws.Cells.Item(46, 22).Value = 1.9   # V46
$ws.Cells.Item(46, 23).Value = 1.9   # W46
$ws.Cells.Item(46, 25).Value = -1   # Y46
$ws.Cells.Item(46, 26).Value = 3.5   # Z46
$ws.Cells.Item(46, 27).Value = -1   # AA46
$ws.Cells.Item(46, 28).Value = 0.925   # AB46
$ws.Cells.Item(46, 29).Value = 0.45   # AC46
$ws.Cells.Item(46, 30).Value = -0.5   # AD46

# Row 58
$ws.Cells.Item(58, 2).Value = 6808899   # B58
$ws.Cells.Item(58, 5).Value = "MKP Pogon Siedlce"   # E58
$ws.Cells.Item(58, 6).Value = "Sandecja Nowy Sacz"   # F58
$ws.Cells.Item(58, 7).Value = 3   # G58
$ws.Cells.Item(58, 8).Value = 0   # H58
$ws.Cells.Item(58, 9).Value = 1   # I58
$ws.Cells.Item(58, 11).Value = "H"   # K58
$ws.Cells.Item(58, 12).Value = 2.25   # L58
$ws.Cells.Item(58, 13).Value = 3.3   # M58
$ws.Cells.Item(58, 14).Value = 2.75   # N58
$ws.Cells.Item(58, 15).Value = 2.25   # O58
$ws.Cells.Item(58, 16).Value = 3.3   # P58
$ws.Cells.Item(58, 17).Value = 2.75   # Q58
$ws.Cells.Item(58, 18).Value = -0.25   # R58
$ws.Cells.Item(58, 19).Value = 2.025   # S58
$ws.Cells.Item(58, 20).Value = 1.775   # T58
$ws.Cells.Item(58, 21).Value = 2.25   # U58
$ws.Cells.Item(58, 22).Value = 1.975   # V58
$ws.Cells.Item(58, 23).Value = 1.825   # W58
$ws.Cells.Item(58, 24).Value = 1.25   # X58
$ws.Cells.Item(58, 26).Value = -1   # Z58
$ws.Cells.Item(58, 27).Value = 1.025   # AA58
$ws.Cells.Item(58, 28).Value = -1   # AB58
$ws.Cells.Item(58, 29).Value = 0.9750000000000001   # AC58
$ws.Cells.Item(58, 30).Value = -1   # AD58

# Row 59
$ws.Cells.Item(59, 2).Value = 7105151   # B59
$ws.Cells.Item(59, 5).Value = "Wisla Pulawy"   # E59
$ws.Cells.Item(59, 6).Value = "Hutnik Krakow"   # F59
$ws.Cells.Item(59, 7).Value = 0   # G59
$ws.Cells.Item(59, 8).Value = 1   # H59
$ws.Cells.Item(59, 9).Value = 0   # I59
$ws.Cells.Item(59, 10).Value = 0   # J59
$ws.Cells.Item(59, 11).Value = "A"   # K59
$ws.Cells.Item(59, 12).Value = 1.8   # L59
$ws.Cells.Item(59, 13).Value = 3.25   # M59
$ws.Cells.Item(59, 14).Value = 4   # N59
$ws.Cells.Item(59, 15).Value = 1.8   # O59
$ws.Cells.Item(59, 16).Value = 3.25   # P59
$ws.Cells.Item(59, 17).Value = 4   # Q59
$ws.Cells.Item(59, 18).Value = -0.5   # R59
$ws.Cells.Item(59, 19).Value = 1.85   # S59
$ws.Cells.Item(59, 20).Value = 1.95   # T59
$ws.Cells.Item(59, 21).Value = 2.5   # U59
$ws.Cells.Item(59, 22).Value = 1.85   # V59
$ws.Cells.Item(59, 23).Value = 1.95   # W59
$ws.Cells.Item(59, 24).Value = -1   # X59
$ws.Cells.Item(59, 26).Value = 3   # Z59
$ws.Cells.Item(59, 27).Value = -1   # AA59
$ws.Cells.Item(59, 28).Value = 0.95   # AB59
$ws.Cells.Item(59, 29).Value = -1   # AC59
$ws.Cells.Item(59, 30).Value = 0.95   # AD59

# Row 60
$ws.Cells.Item(60, 2).Value = 6808625   # B60
$ws.Cells.Item(60, 5).Value = "GKS Jastrzebie"   # E60
$ws.Cells.Item(60, 6).Value = "Polonia Bytom"   # F60
$ws.Cells.Item(60, 7).Value = 4   # G60
$ws.Cells.Item(60, 8).Value = 2   # H60
$ws.Cells.Item(60, 15).Value = 2   # O60
$ws.Cells.Item(60, 16).Value = 3.4   # P60
$ws.Cells.Item(60, 17).Value = 3.1   # Q60
$ws.Cells.Item(60, 19).Value = 1.8   # S60
$ws.Cells.Item(60, 20).Value = 2   # T60
$ws.Cells.Item(60, 21).Value = 2.25   # U60
$ws.Cells.Item(60, 22).Value = 1.875   # V60
$ws.Cells.Item(60, 23).Value = 1.925   # W60
$ws.Cells.Item(60, 24).Value = 1   # X60
$ws.Cells.Item(60, 27).Value = 0.8   # AA60
$ws.Cells.Item(60, 29).Value = 0.875   # AC60

# Row 61
$ws.Cells.Item(61, 2).Value = 6808624   # B61
$ws.Cells.Item(61, 5).Value = "Zaglebie Lubin II"   # E61
$ws.Cells.Item(61, 6).Value = "Olimpia Elblag"   # F61
$ws.Cells.Item(61, 7).Value = 2   # G61
$ws.Cells.Item(61, 8).Value = 1   # H61
$ws.Cells.Item(61, 10).Value = 1   # J61
$ws.Cells.Item(61, 12).Value = 2   # L61
$ws.Cells.Item(61, 13).Value = 3.4   # M61
$ws.Cells.Item(61, 14).Value = 3.1   # N61
$ws.Cells.Item(61, 15).Value = 2.15   # O61
$ws.Cells.Item(61, 16).Value = 3.6   # P61
$ws.Cells.Item(61, 17).Value = 2.7   # Q61
$ws.Cells.Item(61, 19).Value = 1.95   # S61
$ws.Cells.Item(61, 20).Value = 1.85   # T61
$ws.Cells.Item(61, 21).Value = 2.5   # U61
$ws.Cells.Item(61, 22).Value = 1.825   # V61
$ws.Cells.Item(61, 23).Value = 1.975   # W61
$ws.Cells.Item(61, 24).Value = 1.15   # X61
$ws.Cells.Item(61, 27).Value = 0.95   # AA61
$ws.Cells.Item(61, 29).Value = 0.825   # AC61

# Row 65
$ws.Cells.Item(65, 2).Value = 6808631   # B65
$ws.Cells.Item(65, 5).Value = "Olimpia Elblag"   # E65
$ws.Cells.Item(65, 6).Value = "GKS Jastrzebie"   # F65
$ws.Cells.Item(65, 7).Value = 2   # G65
$ws.Cells.Item(65, 8).Value = 1   # H65
$ws.Cells.Item(65, 10).Value = 1   # J65
$ws.Cells.Item(65, 11).Value = "H"   # K65
$ws.Cells.Item(65, 12).Value = 2.1   # L65
$ws.Cells.Item(65, 14).Value = 3.1   # N65
$ws.Cells.Item(65, 15).Value = 2.15   # O65
$ws.Cells.Item(65, 17).Value = 3   # Q65
$ws.Cells.Item(65, 18).Value = -0.25   # R65
$ws.Cells.Item(65, 19).Value = 1.95   # S65
$ws.Cells.Item(65, 20).Value = 1.85   # T65
$ws.Cells.Item(65, 21).Value = 2.25   # U65
$ws.Cells.Item(65, 22).Value = 1.8   # V65
$ws.Cells.Item(65, 23).Value = 2   # W65
$ws.Cells.Item(65, 24).Value = 1.15   # X65
$ws.Cells.Item(65, 26).Value = -1   # Z65
$ws.Cells.Item(65, 27).Value = 0.95   # AA65
$ws.Cells.Item(65, 28).Value = -1   # AB65
$ws.Cells.Item(65, 29).Value = 0.8   # AC65

# Row 66
$ws.Cells.Item(66, 2).Value = 6808898   # B66
$ws.Cells.Item(66, 5).Value = "Sandecja Nowy Sacz"   # E66
$ws.Cells.Item(66, 6).Value = "Wisla Pulawy"   # F66
$ws.Cells.Item(66, 7).Value = 1   # G66
$ws.Cells.Item(66, 8).Value = 2   # H66
$ws.Cells.Item(66, 10).Value = 0   # J66
$ws.Cells.Item(66, 11).Value = "A"   # K66
$ws.Cells.Item(66, 12).Value = 3.1   # L66
$ws.Cells.Item(66, 14).Value = 2.1   # N66
$ws.Cells.Item(66, 15).Value = 2.625   # O66
$ws.Cells.Item(66, 17).Value = 2.4   # Q66
$ws.Cells.Item(66, 18).Value = 0   # R66
$ws.Cells.Item(66, 19).Value = 2   # S66
$ws.Cells.Item(66, 20).Value = 1.8   # T66
$ws.Cells.Item(66, 21).Value = 2.5   # U66
$ws.Cells.Item(66, 22).Value = 2   # V66
$ws.Cells.Item(66, 23).Value = 1.8   # W66
$ws.Cells.Item(66, 24).Value = -1   # X66
$ws.Cells.Item(66, 26).Value = 1.4   # Z66
$ws.Cells.Item(66, 27).Value = -1   # AA66
$ws.Cells.Item(66, 28).Value = 0.8   # AB66
$ws.Cells.Item(66, 29).Value = 1   # AC66

# Row 143
$ws.Cells.Item(143, 2).Value = 6808888   # B143
$ws.Cells.Item(143, 5).Value = "Sandecja Nowy Sacz"   # E143
$ws.Cells.Item(143, 6).Value = "Lech Poznan II"   # F143
$ws.Cells.Item(143, 7).Value = 0   # G143
$ws.Cells.Item(143, 8).Value = 0   # H143
$ws.Cells.Item(143, 9).Value = 0   # I143
$ws.Cells.Item(143, 10).Value = 0   # J143
$ws.Cells.Item(143, 12).Value = 1.85   # L143
$ws.Cells.Item(143, 13).Value = 3.4   # M143
$ws.Cells.Item(143, 14).Value = 3.5   # N143
$ws.Cells.Item(143, 15).Value = 1.85   # O143
$ws.Cells.Item(143, 16).Value = 3.4   # P143
$ws.Cells.Item(143, 17).Value = 3.5   # Q143
$ws.Cells.Item(143, 18).Value = -0.5   # R143
$ws.Cells.Item(143, 19).Value = 1.925   # S143
$ws.Cells.Item(143, 20).Value = 1.875   # T143
$ws.Cells.Item(143, 21).Value = 2.75   # U143
$ws.Cells.Item(143, 22).Value = 1.825   # V143
$ws.Cells.Item(143, 23).Value = 1.975   # W143
$ws.Cells.Item(143, 25).Value = 2.4   # Y143
$ws.Cells.Item(143, 27).Value = -1   # AA143
$ws.Cells.Item(143, 28).Value = 0.875   # AB143
$ws.Cells.Item(143, 29).Value = -1   # AC143
$ws.Cells.Item(143, 30).Value = 0.9750000000000001   # AD143

# Row 144
$ws.Cells.Item(144, 2).Value = 6808197   # B144
$ws.Cells.Item(144, 5).Value = "KKS 1925 Kalisz"   # E144
$ws.Cells.Item(144, 6).Value = "Olimpia Grudziadz"   # F144
$ws.Cells.Item(144, 7).Value = 2   # G144
$ws.Cells.Item(144, 8).Value = 2   # H144
$ws.Cells.Item(144, 9).Value = 2   # I144
$ws.Cells.Item(144, 10).Value = 2   # J144
$ws.Cells.Item(144, 12).Value = 1.666   # L144
$ws.Cells.Item(144, 13).Value = 3.6   # M144
$ws.Cells.Item(144, 14).Value = 4.2   # N144
$ws.Cells.Item(144, 15).Value = 2.3   # O144
$ws.Cells.Item(144, 16).Value = 3.25   # P144
$ws.Cells.Item(144, 17).Value = 2.7   # Q144
$ws.Cells.Item(144, 18).Value = 0   # R144
$ws.Cells.Item(144, 19).Value = 1.775   # S144
$ws.Cells.Item(144, 20).Value = 2.025   # T144
$ws.Cells.Item(144, 21).Value = 2.5   # U144
$ws.Cells.Item(144, 22).Value = 1.85   # V144
$ws.Cells.Item(144, 23).Value = 1.95   # W144
$ws.Cells.Item(144, 25).Value = 2.25   # Y144
$ws.Cells.Item(144, 27).Value = 0   # AA144
$ws.Cells.Item(144, 28).Value = 0   # AB144
$ws.Cells.Item(144, 29).Value = 0.8500000000000001   # AC144
$ws.Cells.Item(144, 30).Value = -1   # AD144

# Row 233
$ws.Cells.Item(233, 2).Value = 6808733   # B233
$ws.Cells.Item(233, 5).Value = "Stomil Olsztyn"   # E233
$ws.Cells.Item(233, 6).Value = "Hutnik Krakow"   # F233
$ws.Cells.Item(233, 8).Value = 0   # H233
$ws.Cells.Item(233, 10).Value = 0   # J233
$ws.Cells.Item(233, 12).Value = 2.375   # L233
$ws.Cells.Item(233, 13).Value = 3.2   # M233
$ws.Cells.Item(233, 14).Value = 2.625   # N233
$ws.Cells.Item(233, 15).Value = 2.4   # O233
$ws.Cells.Item(233, 16).Value = 3.3   # P233
$ws.Cells.Item(233, 17).Value = 2.5   # Q233
$ws.Cells.Item(233, 19).Value = 1.85   # S233
$ws.Cells.Item(233, 20).Value = 1.95   # T233
$ws.Cells.Item(233, 21).Value = 2.25   # U233
$ws.Cells.Item(233, 22).Value = 2.025   # V233
$ws.Cells.Item(233, 23).Value = 1.775   # W233
$ws.Cells.Item(233, 24).Value = 1.4   # X233
$ws.Cells.Item(233, 27).Value = 0.8500000000000001   # AA233
$ws.Cells.Item(233, 29).Value = 1.025   # AC233

# Row 234
$ws.Cells.Item(234, 2).Value = 6808048   # B234
$ws.Cells.Item(234, 5).Value = "Radunia Stezyca"   # E234
$ws.Cells.Item(234, 6).Value = "Skra Czestochowa"   # F234
$ws.Cells.Item(234, 8).Value = 2   # H234
$ws.Cells.Item(234, 10).Value = 1   # J234
$ws.Cells.Item(234, 12).Value = 2   # L234
$ws.Cells.Item(234, 13).Value = 3.1   # M234
$ws.Cells.Item(234, 14).Value = 3.4   # N234
$ws.Cells.Item(234, 15).Value = 2.45   # O234
$ws.Cells.Item(234, 16).Value = 3   # P234
$ws.Cells.Item(234, 17).Value = 2.875   # Q234
$ws.Cells.Item(234, 19).Value = 1.75   # S234
$ws.Cells.Item(234, 20).Value = 2.05   # T234
$ws.Cells.Item(234, 21).Value = 2   # U234
$ws.Cells.Item(234, 22).Value = 1.9   # V234
$ws.Cells.Item(234, 23).Value = 1.9   # W234
$ws.Cells.Item(234, 24).Value = 1.45   # X234
$ws.Cells.Item(234, 27).Value = 0.75   # AA234
$ws.Cells.Item(234, 29).Value = 0.8999999999999999   # AC234
